# Fill in the missing PriceChange/UpDown values for row 6, and append a new
# data row (row 7) for 2016-10-05, matching a fresh scan run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: complete the trailing PriceChange / UpDown columns ---
$ws.Range("X6").Value = 0.059999000000001246
$ws.Range("Y6").Value = "Up"

# --- Row 7: brand new scan result row ---
$ws.Range("A7").Value = 42648.886562500003
$ws.Range("B7").Value = 11
$ws.Range("C7").Value = "Buy"
$ws.Range("D7").Value = 44
$ws.Range("E7").Value = 10294
$ws.Range("F7").Value = 1667
$ws.Range("G7").Value = 68
$ws.Range("H7").Value = 29
$ws.Range("I7").Value = 92
$ws.Range("J7").Value = 7
$ws.Range("K7").Value = 18929
$ws.Range("L7").Value = 299
$ws.Range("M7").Value = 130
$ws.Range("N7").Value = 93
$ws.Range("O7").Value = 8
$ws.Range("P7").Value = "Noun"
$ws.Range("Q7").Value = 41.162214763508182
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0.0616
$ws.Range("T7").Value = -0.032
$ws.Range("U7").Value = 2.2599999999999998
$ws.Range("V7").Value = "N/A"
$ws.Range("W7").Value = 0

# Carry the existing number formats down from row 2 so no new style entries
# are minted (matches the date-serial style on column A and the percentage
# style on columns S:T used throughout the rest of the table).
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)

$ws.Range("S2:T2").Copy()
$ws.Range("S7:T7").PasteSpecial(-4122)
